$d = $word.ActiveDocument

# --- 1. Change the resistor value "10k" -> "14.7k" -------------------------
# ("1 - 10k Resistor (RB)" becomes "1 - 14.7k Resistor (RB)")
$hit = $d.Content.Duplicate
$found = $hit.Find.Execute("10k", $true, $false, $false, $false, $false, `
                            $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find '10k' to edit"
}

# $hit spans exactly "10k" (3 chars): the "1" run, the "0" run, the "k" run.
$zeroRun = $d.Range($hit.Start + 1, $hit.Start + 2)
if ($zeroRun.Text -ne "0") {
    throw "Unexpected text at target position: [$($zeroRun.Text)]"
}
$zeroRun.Text = "4.7"

# Editing the "0" run merges it together with the neighboring "1" and "k"
# text into a single run ("14.7k"). Restore the original run boundary between
# "1" and "4.7" by bookmarking (and then un-bookmarking) the split point --
# inserting/removing a bookmark there forces the text back apart into two
# runs without leaving any trace of the temporary bookmark behind.
$splitPos = $d.Range($hit.Start + 1, $hit.Start + 1)
$d.Bookmarks.Add("zzTempSplit", $splitPos) | Out-Null
$d.Bookmarks("zzTempSplit").Delete()

# --- 2. Move the "_GoBack" bookmark ----------------------------------------
# It now needs to sit between "4.7" and "k" (i.e. "1" / "4.7" / [bookmark] / "k"),
# rather than at the end of the "Note: ..." paragraph. Re-locate "14.7k" since
# the run layout changed above.
$hit2 = $d.Content.Duplicate
$found2 = $hit2.Find.Execute("14.7k", $true, $false, $false, $false, $false, `
                              $true, 1, $false, "", 0)
if (-not $found2) {
    throw "Could not find '14.7k' after editing"
}
$bookmarkPos = $d.Range($hit2.Start + 4, $hit2.Start + 4)

# Adding a bookmark named "_GoBack" here both splits the "4.7k" run into
# "4.7" and "k", and removes the old "_GoBack" bookmark from its previous
# location (bookmark names must be unique in the document).
$d.Bookmarks.Add("_GoBack", $bookmarkPos) | Out-Null

Write-Output "Resistor value updated and _GoBack bookmark relocated."
